$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (existing row): ticks 20 -> 8, expected (tick size * ticks) -> 700
$ws.Range("A9").Value = 8
$ws.Range("B9").Formula = '=$A$4*A9'
$ws.Range("C9").Value = 700

# New row 10
$ws.Range("A10").Value = 13
$ws.Range("B10").Formula = '=$A$4*A10'
$ws.Range("C10").Value = 1100

# New row 11
$ws.Range("A11").Value = 20
$ws.Range("B11").Formula = '=$A$4*A11'
$ws.Range("C11").Value = 1700

# New row 12
$ws.Range("A12").Value = 50
$ws.Range("B12").Formula = '=$A$4*A12'
$ws.Range("C12").Value = 4200

# New row 13
$ws.Range("A13").Value = 200
$ws.Range("B13").Formula = '=$A$4*A13'
$ws.Range("C13").Value = 17000
